# Fix the "Conpany & APR" typo to "Company & APR" in the table header cell.
# The canonical edit (per the target diff) splits the single run into two
# runs - "Com" and "pany & APR" - both keeping the original bold / sz=24
# run formatting, while leaving the paragraph's own properties untouched.

$d = $word.ActiveDocument

# Locate the exact run text and narrow $r down to just that span.
$r = $d.Content
$r.Find.Execute("Conpany & APR") | Out-Null

# Rebuild that paragraph's content as two runs with identical formatting,
# matching the way the paragraph already looked (same w:p attributes and
# w:pPr), but with the text split after "Com".
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00067FD9" w:rsidRPr="0076121A" w:rsidRDefault="00067FD9" w:rsidP="0076121A"><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/></w:rPr><w:t>Com</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="24"/></w:rPr><w:t>pany &amp; APR</w:t></w:r></w:p>
'@

$r.InsertXML($xml) | Out-Null
